$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer readings (x, y, z) captured on May 9th, to be inserted
# directly after the header row, pushing all existing data rows down by 10.
$data = @(
    @(-0.3973007202148437, 1.798293828964233, 1.168385148048401),
    @(-0.8141142129898071, 1.826734185218811, 1.343665383756161),
    @(-0.9248467683792115, 1.964664489030838, 1.209049716591835),
    @(-0.7549184560775757, 1.929333925247193, 1.357575602829456),
    @(-0.4957029819488525, 1.878820419311524, 1.164325326681137),
    @(-0.8541634678840635, 1.869282335042953, 1.307794235646725),
    @(-0.8955824375152587, 1.77057421207428, 1.258034527301789),
    @(-0.5728458166122438, 1.762963086366654, 1.322432711720467),
    @(-0.5884580612182617, 1.777032017707825, 1.320214748382568),
    @(-0.7303044199943542, 1.775961980223656, 1.753339484333992)
)

$rowCount = $data.Count

# Insert blank rows right below the header (row 1), shifting existing data down.
$insertRange = $ws.Range("A2:A$(1 + $rowCount)")
$insertRange.EntireRow.Insert()

# The inserted rows pick up formatting from the row above (the header); strip
# it back off so the new data rows stay unstyled, like the other data rows.
$ws.Range("A2:C$(1 + $rowCount)").ClearFormats()

# Fill in the newly inserted rows with the new readings.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
